$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the two worker data rows (row 16 and row 17): previous database
# entries are removed and the new ones added, swapping positions so the
# "base de datos" now lists YAIR FERNEL OSORIO RAMIREZ first, then
# DANIELA MARTINEZ ALVAREZ.

$ws.Range("C16").Value = "1143131446"
$ws.Range("D16").Value = "YAIR FERNEL OSORIO RAMIREZ"
$ws.Range("E16").Value = "2412"
$ws.Range("F16").Value = 388950
$ws.Range("G16").Value = 9723763

$ws.Range("C17").Value = "1007588337"
$ws.Range("D17").Value = "DANIELA MARTINEZ ALVAREZ"
$ws.Range("E17").Value = "2303"
$ws.Range("F17").Value = 7208
$ws.Range("G17").Value = 6394174
